$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.984.71'
$ws.Range("E2").Value = '  +1.49%  '
$ws.Range("D3").Value = '3.121.90'
$ws.Range("E3").Value = '  +0.36%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '586.18'
$ws.Range("D6").Value = '146.40'
$ws.Range("E6").Value = '  +1.65%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '3.113.79'
$ws.Range("E8").Value = '  +0.38%  '
$ws.Range("D9").Value = '0.530'
$ws.Range("E9").Value = '  -0.35%  '
$ws.Range("E10").Value = '  +9.42%  '
$ws.Range("E11").Value = '  +1.18%  '
$ws.Range("E12").Value = '  -0.96%  '
$ws.Range("E13").Value = '  +2.51%  '
$ws.Range("D14").Value = '37.31'
$ws.Range("E14").Value = '  +4.56%  '
$ws.Range("E15").Value = '  -0.77%  '
$ws.Range("D16").Value = '3.638.56'
$ws.Range("E16").Value = '  +0.34%  '
$ws.Range("D17").Value = '63.859.12'
$ws.Range("E17").Value = '  +1.40%  '
$ws.Range("D18").Value = '7.13'
$ws.Range("E18").Value = '  -1.88%  '
$ws.Range("D19").Value = '3.119.28'
$ws.Range("E19").Value = '  +0.39%  '
$ws.Range("D20").Value = '463.65'
$ws.Range("E20").Value = '  +2.07%  '
$ws.Range("D21").Value = '14.31'
$ws.Range("E21").Value = '  +1.33%  '
$ws.Range("E22").Value = '  -0.55%  '
$ws.Range("E23").Value = '  -0.87%  '
$ws.Range("E24").Value = '  -3.32%  '
$ws.Range("D25").Value = '81.71'
$ws.Range("E25").Value = '  -0.53%  '
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("D27").Value = '8.94'
$ws.Range("E27").Value = '  +8.35%  '
$ws.Range("E28").Value = '  -0.51%  '
$ws.Range("E29").Value = '  -1.51%  '
$ws.Range("E30").Value = '  -0.01%  '
$ws.Range("E31").Value = '  +0.22%  '
$ws.Range("E32").Value = '  -0.17%  '
$ws.Range("E33").Value = '  -2.78%  '
$ws.Range("D34").Value = '0.0₃0866'
$ws.Range("E34").Value = '  +6.63%  '
$ws.Range("E35").Value = '  -1.14%  '
$ws.Range("E36").Value = '  +1.32%  '
$ws.Range("D37").Value = '3.39'
$ws.Range("E37").Value = '  +9.44%  '
$ws.Range("D38").Value = '6.05'
$ws.Range("E38").Value = '  -0.12%  '
$ws.Range("E39").Value = '  -0.24%  '
$ws.Range("D40").Value = '446.41'
$ws.Range("E40").Value = '  +4.53%  '
$ws.Range("E41").Value = '  -0.81%  '
$ws.Range("D42").Value = '0.0371'
$ws.Range("E42").Value = '  -0.43%  '
$ws.Range("D43").Value = '2.874.75'
$ws.Range("E43").Value = '  -3.19%  '
$ws.Range("E44").Value = '  -0.08%  '
$ws.Range("E45").Value = '  -1.21%  '
$ws.Range("E46").Value = '  -0.69%  '
$ws.Range("D47").Value = '35.75'
$ws.Range("E47").Value = '  +3.23%  '
$ws.Range("E48").Value = '  +0.04%  '
$ws.Range("D49").Value = '123.33'
$ws.Range("E49").Value = '  -1.19%  '
$ws.Range("E50").Value = '  -0.73%  '
$ws.Range("D51").Value = '24.64'
$ws.Range("E51").Value = '  -1.45%  '
